$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Joins": remove the REVISIONS/DOCUMENTS join row (row 2), and
# mark the remaining join (DOCUMENTS/UNSTRUCTURED_DATA) as not joined.
# ---------------------------------------------------------------------
$wsJoins = $wb.Worksheets.Item("Joins")

$wsJoins.Rows.Item(2).Delete()
$wsJoins.Range("E2").Value = "n"

[void]$wsJoins.Range("E2").Select()

# ---------------------------------------------------------------------
# Sheet "Input": insert a new "Label" column (D) between "Input" (C) and
# "Date" (old D, now E). Add a new row 3 with a second document-id/name
# pair.
# ---------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

$wsInput.Columns.Item(4).Insert()
$wsInput.Columns.Item(4).ColumnWidth = $wsInput.Columns.Item(3).ColumnWidth

$wsInput.Range("D1").Value = "Label"
$wsInput.Range("D2").Value = "Document Id"
$wsInput.Range("C3").Value = "DORIGINALNAME"
$wsInput.Range("D3").Value = "Document Name"

[void]$wsInput.Range("D2:D3").Select()

# ---------------------------------------------------------------------
# Sheet "Output": insert a new "Label" column (C) between "Output" (B)
# and "Download" (old C, now D). Add matching Document Id/Name values.
# ---------------------------------------------------------------------
$wsOutput = $wb.Worksheets.Item("Output")

$wsOutput.Columns.Item(3).Insert()
$wsOutput.Columns.Item(3).ColumnWidth = $wsOutput.Columns.Item(2).ColumnWidth

$wsOutput.Range("C1").Value = "Label"
$wsOutput.Range("C2").Value = "Document Id"
$wsOutput.Range("C3").Value = "Document Name"

[void]$wsOutput.Range("C2:C3").Select()

# ---------------------------------------------------------------------
# Make "Output" the active sheet/tab (was "Joins").
# ---------------------------------------------------------------------
[void]$wsOutput.Activate()
